# LQA_Tester_ProgressTracker.xlsx refresh
#
# Recomputed stats after clearing any AutoFilter on the master data sheet
# before re-sorting (Item category A-Z) and rebuilding DAILY / TOTAL.
# Writes the refreshed numbers for the latest data pull (01/09) straight
# into the three sheets: DAILY, TOTAL, and the hidden _DAILY_DATA source.

$wb = $excel.ActiveWorkbook

# --- DAILY: row 8 (01/09) now has real numbers for Doni, Eric, Lisa, and
#     updated totals for John/Mike/Paul/Pending instead of placeholder "--"
$ws = $wb.Worksheets.Item("DAILY")
$ws.Range("F8").Value = 385   # Doni - Done
$ws.Range("G8").Value = 104   # Doni - Issues
$ws.Range("H8").Value = 3     # Eric - Done
$ws.Range("I8").Value = 1     # Eric - Issues
$ws.Range("J8").Value = 350   # John - Done
$ws.Range("K8").Value = 5     # John - Issues
$ws.Range("N8").Value = 37    # Lisa - Done
$ws.Range("O8").Value = 17    # Lisa - Issues
$ws.Range("P8").Value = 253   # Mike - Done
$ws.Range("Q8").Value = 86    # Mike - Issues
$ws.Range("R8").Value = 254   # Paul - Done
$ws.Range("X8").Value = 280   # Pending

# --- TOTAL: per-tester completion/issue stats plus SUBTOTAL/TOTAL rollups
$ws = $wb.Worksheets.Item("TOTAL")

# Doni
$ws.Range("B5").Value = 99.7
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = 385
$ws.Range("E5").Value = 104
$ws.Range("F5").Value = 276
$ws.Range("G5").Value = 5
$ws.Range("K5").Value = 104

# Eric
$ws.Range("B6").Value = 37.5
$ws.Range("C6").Value = 100
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("K6").Value = 1

# John
$ws.Range("D7").Value = 354
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 324
$ws.Range("G7").Value = 23
$ws.Range("K7").Value = 7

# Lisa
$ws.Range("B9").Value = 97.4
$ws.Range("C9").Value = 100
$ws.Range("D9").Value = 37
$ws.Range("E9").Value = 17
$ws.Range("F9").Value = 20
$ws.Range("K9").Value = 17

# Mike
$ws.Range("B10").Value = 96.6
$ws.Range("D10").Value = 253
$ws.Range("E10").Value = 86
$ws.Range("F10").Value = 167
$ws.Range("K10").Value = 86

# Paul
$ws.Range("B11").Value = 96.9
$ws.Range("D11").Value = 254
$ws.Range("F11").Value = 191

# SUBTOTAL
$ws.Range("B12").Value = 95.1
$ws.Range("D12").Value = 1326
$ws.Range("E12").Value = 285
$ws.Range("F12").Value = 1012
$ws.Range("G12").Value = 29
$ws.Range("K12").Value = 285

# GRAND TOTAL
$ws.Range("B15").Value = 95.1
$ws.Range("D15").Value = 1326
$ws.Range("E15").Value = 285
$ws.Range("F15").Value = 1012
$ws.Range("G15").Value = 29
$ws.Range("K15").Value = 285

# --- _DAILY_DATA (hidden master): per-(date,user,category) rows that feed
#     the DAILY/TOTAL rollups above
$ws = $wb.Worksheets.Item("_DAILY_DATA")

# Row 6: Doni / Knowledge (01/09)
$ws.Range("D6").Value = 386
$ws.Range("E6").Value = 385
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 276
$ws.Range("H6").Value = 5

# Row 8: Mike / Region (01/09)
$ws.Range("E8").Value = 253
$ws.Range("F8").Value = 86
$ws.Range("G8").Value = 167

# Row 9: Lisa / Region (01/09)
$ws.Range("E9").Value = 37
$ws.Range("F9").Value = 17
$ws.Range("G9").Value = 20

# Row 10: Paul / Region (01/09)
$ws.Range("E10").Value = 254
$ws.Range("G10").Value = 191

# Row 12: Eric / Quest (01/09)
$ws.Range("D12").Value = 8
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1

# Row 13: John / Quest (01/09)
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 4
